# The survey data included 4 rows (the original rows 912-915) for owls
# recorded on Route/Station 99. Per the commit message ("updated without
# Route 99 owl"), those rows are removed entirely, and all following rows
# shift up by 4 (e.g. old row 916 -> new row 912, ... old row 944 -> new
# row 940). The table's defined name and the sheet's used-range dimension
# both need to reflect the new last row (940 instead of 944).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 4 rows belonging to Route 99 (original rows 912-915).
# This shifts every subsequent row up by 4 and also updates the sheet's
# <dimension> automatically.
$ws.Rows.Item(912).Resize(4).Delete()

# Update the workbook-level defined name "Owls_Table" so it refers to the
# new, smaller range (it previously referred to $A$1:$J$944).
$wb.Names.Item("Owls_Table").RefersTo = "='Owls_Table'!`$A`$1:`$J`$940"
